$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 450
$ws.Range("F4").Value = 128
$ws.Range("F5").Value = 366
$ws.Range("F6").Value = 587
$ws.Range("F10").Value = 399
$ws.Range("F12").Value = 771
$ws.Range("F13").Value = 782
$ws.Range("F16").Value = 1540
$ws.Range("F17").Value = 1540
$ws.Range("F18").Value = 946
$ws.Range("F22").Value = 360
$ws.Range("F25").Value = 113
$ws.Range("F26").Value = 6739
$ws.Range("F27").Value = 5141
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 458
$ws.Range("F29").Value = 149
$ws.Range("F30").Value = 489
$ws.Range("F31").Value = 213
$ws.Range("F35").Value = 30
$ws.Range("F37").Value = 1311
$ws.Range("F38").Value = 200
$ws.Range("F39").Value = 259
$ws.Range("F40").Value = 630
$ws.Range("F43").Value = 268
$ws.Range("F45").Value = 156
$ws.Range("F46").Value = 67

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 14
$ws.Range("F6").Value = 42
$ws.Range("F18").Value = 253

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 172
$ws.Range("F3").Value = 2479
$ws.Range("F4").Value = 214
$ws.Range("F5").Value = 81

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 450
$ws.Range("F6").Value = 214
$ws.Range("F7").Value = 81
$ws.Range("F8").Value = 366
$ws.Range("F9").Value = 587
$ws.Range("F13").Value = 14
$ws.Range("F14").Value = 399
$ws.Range("F16").Value = 772
$ws.Range("F17").Value = 782
$ws.Range("F20").Value = 1540
$ws.Range("F21").Value = 1540
$ws.Range("F22").Value = 946
$ws.Range("F24").Value = 360
$ws.Range("F26").Value = 113
$ws.Range("F27").Value = 42
$ws.Range("F29").Value = 6740
$ws.Range("F30").Value = 5142
$ws.Range("F32").Value = 30
$ws.Range("F33").Value = 1311
$ws.Range("F34").Value = 200
$ws.Range("F36").Value = 259
$ws.Range("F38").Value = 630
$ws.Range("F43").Value = 268
$ws.Range("F44").Value = 156
$ws.Range("F45").Value = 67
$ws.Range("F49").Value = 253
